# Update the account statistics table in Accounts.docx.
# Table layout (row 1 = header):
#   Col 1: ID              Col 2: Tên đăng nhập     Col 3: Mật khẩu
#   Col 4: Số tiền         Col 5: Số trận đã chơi   Col 6: Số trận thắng
#   Col 7: Số trận thua    Col 8: Tổng số tiền thắng cược
#   Col 9: Tổng số tiền đã thua

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Row 2 (account 001)
$t.Cell(2, 4).Range.Text = "140"
$t.Cell(2, 5).Range.Text = "30"
$t.Cell(2, 7).Range.Text = "13"
$t.Cell(2, 8).Range.Text = "667784"
$t.Cell(2, 9).Range.Text = "676174"

# Row 3 (account 002)
$t.Cell(3, 4).Range.Text = "53"
$t.Cell(3, 5).Range.Text = "10"
$t.Cell(3, 7).Range.Text = "5"
$t.Cell(3, 9).Range.Text = "1543"
